# Change the table style used by the three tables on slides 14-16
# from "{0D6FE889-D6C2-429E-AA0F-6331EC94F161}" (Table_0) to the
# built-in style "{EC30C1BF-455E-41C5-87D3-AAA1B03F345C}".
$p = $ppt.ActivePresentation

$oldStyleId = "{0D6FE889-D6C2-429E-AA0F-6331EC94F161}"
$newStyleId = "{EC30C1BF-455E-41C5-87D3-AAA1B03F345C}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# Swap the presentation's design theme colour scheme ("Integral" / "Red
# Violet") for the stock "Office Theme" colour scheme that used to live
# only on the notes master theme (fonts/effects are identical between the
# two themes already, only the 12 theme colours differ).
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeElements.ThemeColorScheme

# RGB() values below use the same BGR-packed-integer encoding the
# PowerPoint object model uses for ColorFormat.RGB (0x00BBGGRR).
$colorScheme.Item(1).RGB  = 0        # dk1      -> 000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      -> FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      -> 44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      -> E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  -> 5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  -> ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  -> A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  -> FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  -> 4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  -> 70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    -> 0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink -> 954F72
